# Apply updated crypto price / volume data as published by the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.281.08'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '1.807.04'
$ws.Range("E3").Value = '  +3.39%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''337.35'
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '''0.4669'
$ws.Range("E7").Value = '  +20.90%  '
$ws.Range("D8").Value = '''0.3861'
$ws.Range("E8").Value = '  +14.19%  '
$ws.Range("D9").Value = '''45.35'
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").Value = '''1.159'
$ws.Range("E10").Value = '  +4.23%  '
$ws.Range("D11").Value = '''0.07603'
$ws.Range("E11").Value = '  +5.79%  '
$ws.Range("D12").Value = '''22.52'
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("D13").Value = '''1.003'
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("D14").Value = '''6.360'
$ws.Range("E14").Value = '  +3.03%  '
$ws.Range("D15").Value = '''7.499'
$ws.Range("E15").Value = '  +5.98%  '
$ws.Range("D16").Value = '1.807.66'
$ws.Range("E16").Value = '  +3.37%  '
$ws.Range("E17").Value = '  +3.75%  '
$ws.Range("D18").Value = '''0.06735'
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("D19").Value = '''81.79'
$ws.Range("E19").Value = '  +3.31%  '
$ws.Range("D20").Value = '''0.9998'
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '''17.61'
$ws.Range("E21").Value = '  +5.19%  '
$ws.Range("D22").Value = '''6.437'
$ws.Range("E22").Value = '  +4.41%  '
$ws.Range("D23").Value = '28.269.50'
$ws.Range("E23").Value = '  +1.35%  '
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("D25").Value = '''2.430'
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("D26").Value = '''20.75'
$ws.Range("E26").Value = '  +4.71%  '
$ws.Range("D27").Value = '''153.28'
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("D28").Value = '''2.377'
$ws.Range("E28").Value = '  +3.64%  '
$ws.Range("D29").Value = '2.012.65'
$ws.Range("E29").Value = '  +3.30%  '
$ws.Range("D30").Value = '''133.10'
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("D31").Value = '''1.257'
$ws.Range("E31").Value = '  -2.38%  '
$ws.Range("D32").Value = '''4.038'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("D33").Value = '''0.09664'
$ws.Range("E33").Value = '  +10.07%  '
$ws.Range("D34").Value = '''5.865'
$ws.Range("E34").Value = '  +0.90%  '
$ws.Range("D35").Value = '''0.2383'
$ws.Range("E35").Value = '  +14.04%  '
$ws.Range("D36").Value = '''0.06365'
$ws.Range("E36").Value = '  +4.08%  '
$ws.Range("D37").Value = '''12.11'
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '''5.285'
$ws.Range("E38").Value = '  +3.20%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.02354'
$ws.Range("E39").Value = '  +3.54%  '
$ws.Range("D40").Value = '''0.6627'
$ws.Range("E40").Value = '  +1.67%  '
$ws.Range("D41").Value = '''1.235'
$ws.Range("E41").Value = '  +2.74%  '
$ws.Range("D42").Value = '''8.421'
$ws.Range("E42").Value = '  +5.21%  '
$ws.Range("E43").Value = '  -2.92%  '
$ws.Range("D44").Value = '''14.30'
$ws.Range("E44").Value = '  +4.24%  '
$ws.Range("D45").Value = '''1.000'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").Value = '''0.6144'
$ws.Range("E46").Value = '  +2.01%  '
$ws.Range("E47").Value = '  +0.99%  '
$ws.Range("D48").Value = '''131.41'
$ws.Range("E48").Value = '  +3.80%  '
$ws.Range("D49").Value = '''2.045'
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("D50").Value = '''1.183'
$ws.Range("E50").Value = '  +1.74%  '
$ws.Range("D51").Value = '''0.07141'
$ws.Range("E51").Value = '  +2.61%  '
